{"js": "// Replace legacy \"$Var\" placeholders with Jinja2-style \"{{ Var }}\"\n// placeholders (docxtpl syntax), throughout the document body.\nconst varNames = [\n  \"EmpName\",\n  \"DOB\",\n  \"EmpID\",\n  \"Department\",\n  \"Position\",\n  \"Salary\",\n  \"JoinDate\",\n  \"Manager\",\n];\n\nconst body = context.document.body;\n\nfor (const name of varNames) {\n  const searchText = \"$\" + name;\n  const replacement = \"{{ \" + name + \" }}\";\n\n  const results = body.search(searchText, {\n    matchCase: true,\n    matchWholeWord: false,\n  });\n  results.load(\"items\");\n  await context.sync();\n\n  for (let i = 0; i < results.items.length; i++) {\n    results.items[i].insertText(replacement, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "# Replace legacy \"$Var\" placeholders with Jinja2-style \"{{ Var }}\"\n# placeholders (docxtpl syntax), throughout the document.\n$d = $word.ActiveDocument\n\n$varNames = @(\n    'EmpName',\n    'DOB',\n    'EmpID',\n    'Department',\n    'Position',\n    'Salary',\n    'JoinDate',\n    'Manager'\n)\n\nforeach ($name in $varNames) {\n    $searchText = '$' + $name\n    $replaceText = '{{ ' + $name + ' }}'\n\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Text = $searchText\n    $find.Replacement.Text = $replaceText\n    $find.Forward = $true\n    $find.Wrap = 1\n    $find.Format = $false\n    $find.MatchCase = $true\n    $find.MatchWholeWord = $false\n    $find.MatchWildcards = $false\n    $find.MatchSoundsLike = $false\n    $find.MatchAllWordForms = $false\n\n    $find.Execute($find.Text, $find.MatchCase, $find.MatchWholeWord, $find.MatchWildcards, $find.MatchSoundsLike, $find.MatchAllWordForms, $find.Forward, $find.Wrap, $find.Format, $find.Replacement.Text, 2)\n}\n"}
